$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (values such as "579.40" or "0.0360" must not be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.828.13"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "3.089.58"
$ws.Range("E3").Value = "  +5.09%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "579.40"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "167.77"
$ws.Range("E6").Value = "  +5.70%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.084.24"
$ws.Range("E8").Value = "  +4.99%  "
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "36.37"
$ws.Range("E14").Value = "  +5.93%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "3.602.88"
$ws.Range("E16").Value = "  +5.09%  "
$ws.Range("D17").Value = "66.848.89"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "3.091.14"
$ws.Range("E19").Value = "  +5.22%  "
$ws.Range("D20").Value = "16.13"
$ws.Range("E20").Value = "  +4.13%  "
$ws.Range("D21").Value = "467.01"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("E22").Value = "  +3.81%  "
$ws.Range("D23").Value = "7.51"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").Value = "83.78"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  +6.36%  "
$ws.Range("D26").Value = "13.07"
$ws.Range("E26").Value = "  +7.88%  "
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  +3.88%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("D37").Value = "5.89"
$ws.Range("E37").Value = "  +2.67%  "

# Rows 38 and 39 swap: Arweave moves up to rank 38, Stacks moves to rank 39
$ws.Range("B38").Value = "Arweave"
$ws.Range("C38").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D38").Value = "46.97"
$ws.Range("E38").Value = "  +3.36%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.11"
$ws.Range("E39").Value = "  +6.32%  "

$ws.Range("D40").Value = "50.34"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("E41").Value = "  +6.10%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "0.0360"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").Value = "382.27"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "2.778.90"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").Value = "135.19"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "24.97"
$ws.Range("E50").Value = "  +6.39%  "
$ws.Range("E51").Value = "  +1.79%  "
